$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 524.6429000000001
$ws.Range("I12").Value = 367.45456
$ws.Range("K12").Value = 367.45456
$ws.Range("M12").Value = -197.45456
$ws.Range("H64").Value = 4659.7
$ws.Range("I64").Value = 3599.5
$ws.Range("J64").Value = 6250
$ws.Range("K64").Value = 3599.5
$ws.Range("L64").Value = 6250
$ws.Range("M64").Value = -3351.5
$ws.Range("N64").Value = -6746
$ws.Range("H67").Value = 4659.7
$ws.Range("I67").Value = 3599.5
$ws.Range("J67").Value = 6250
$ws.Range("K67").Value = 3599.5
$ws.Range("L67").Value = 6250
$ws.Range("M67").Value = -2741.5
$ws.Range("N67").Value = -7966
$ws.Range("H98").Value = 1676.7407
$ws.Range("I98").Value = 1329.579
$ws.Range("J98").Value = 2501.25
$ws.Range("K98").Value = 1329.579
$ws.Range("L98").Value = 2501.25
$ws.Range("M98").Value = 168.421
$ws.Range("N98").Value = -5497.25
$ws.Range("H110").Value = 55000
$ws.Range("J110").Value = 55000
$ws.Range("L110").Value = 55000
$ws.Range("N110").Value = -63180
$ws.Range("H120").Value = 170000
$ws.Range("J120").Value = 170000
$ws.Range("L120").Value = 170000
$ws.Range("N120").Value = -179676
$ws.Range("H122").Value = 1676.7407
$ws.Range("I122").Value = 1329.579
$ws.Range("J122").Value = 2501.25
$ws.Range("K122").Value = 3988.737
$ws.Range("L122").Value = 7503.75
$ws.Range("M122").Value = -1538.737
$ws.Range("N122").Value = -12403.75
$ws.Range("H133").Value = 81012.5
$ws.Range("J133").Value = 81012.5
$ws.Range("L133").Value = 81012.5
$ws.Range("N133").Value = -91132.5
$ws.Range("H135").Value = 649.5
$ws.Range("I135").Value = 649.5
$ws.Range("K135").Value = 5845.5
$ws.Range("M135").Value = -3310.5
$ws.Range("H137").Value = 3438.04
$ws.Range("I137").Value = 3639.1592
$ws.Range("J137").Value = 1963.1666
$ws.Range("K137").Value = 10917.4776
$ws.Range("L137").Value = 5889.4998
$ws.Range("M137").Value = -8367.4776
$ws.Range("N137").Value = -10989.4998
$ws.Range("H138").Value = 13702405
$ws.Range("J138").Value = 4016.261
$ws.Range("L138").Value = 12048.783
$ws.Range("N138").Value = -22328.783

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1408.11
$ws.Range("I32").Value = 1408.11
$ws.Range("K32").Value = 1408.11
$ws.Range("M32").Value = -1121.11
$ws.Range("H131").Value = 55950
$ws.Range("J131").Value = 55950
$ws.Range("L131").Value = 55950
$ws.Range("N131").Value = -66030
$ws.Range("H132").Value = 1998.1025
$ws.Range("I132").Value = 1943.0555
$ws.Range("K132").Value = 5829.166499999999
$ws.Range("M132").Value = -3299.166499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 14998
$ws.Range("I26").Value = 14998
$ws.Range("K26").Value = 14998
$ws.Range("M26").Value = -14706
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
$ws.Range("H134").Value = 2119.1633
$ws.Range("I134").Value = 2119.1633
$ws.Range("K134").Value = 6357.4899
$ws.Range("M134").Value = -3822.4899

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2417.3333
$ws.Range("I58").Value = 1735.0714
$ws.Range("J58").Value = 4805.25
$ws.Range("K58").Value = 1735.0714
$ws.Range("L58").Value = 4805.25
$ws.Range("M58").Value = -1532.0714
$ws.Range("N58").Value = -5211.25
$ws.Range("H99").Value = 11402.839
$ws.Range("I99").Value = 8561.375
$ws.Range("K99").Value = 8561.375
$ws.Range("M99").Value = -7063.375
$ws.Range("H126").Value = 11402.839
$ws.Range("I126").Value = 8561.375
$ws.Range("K126").Value = 25684.125
$ws.Range("M126").Value = -23214.125
$ws.Range("H132").Value = 3409.3635
$ws.Range("I132").Value = 3216.39
$ws.Range("J132").Value = 5035.857
$ws.Range("K132").Value = 9649.17
$ws.Range("L132").Value = 15107.571
$ws.Range("M132").Value = -7119.17
$ws.Range("N132").Value = -20167.571
$ws.Range("H136").Value = 2417.3333
$ws.Range("I136").Value = 1735.0714
$ws.Range("J136").Value = 4805.25
$ws.Range("K136").Value = 5205.2142
$ws.Range("L136").Value = 14415.75
$ws.Range("M136").Value = -2655.2142
$ws.Range("N136").Value = -19515.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 4800
$ws.Range("I64").Value = 4800
$ws.Range("K64").Value = 14400
$ws.Range("M64").Value = -14130
$ws.Range("H67").Value = 4800
$ws.Range("I67").Value = 4800
$ws.Range("K67").Value = 14400
$ws.Range("M67").Value = -13464
$ws.Range("H87").Value = 6999
$ws.Range("I87").Value = 6999
$ws.Range("K87").Value = 20997
$ws.Range("M87").Value = -19749
$ws.Range("H90").Value = 6999
$ws.Range("I90").Value = 6999
$ws.Range("K90").Value = 62991
$ws.Range("M90").Value = -56751
$ws.Range("H114").Value = 1147.75
$ws.Range("J114").Value = 2028.625
$ws.Range("L114").Value = 6085.875
$ws.Range("N114").Value = -12593.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 54750
$ws.Range("J46").Value = 94500
$ws.Range("L46").Value = 94500
$ws.Range("N46").Value = -94812
$ws.Range("H47").Value = 17290
$ws.Range("I47").Value = 15000
$ws.Range("K47").Value = 15000
$ws.Range("M47").Value = -14432
$ws.Range("H107").Value = 965.25
$ws.Range("I107").Value = 865.3333
$ws.Range("J107").Value = 1025.2
$ws.Range("K107").Value = 865.3333
$ws.Range("L107").Value = 1025.2
$ws.Range("M107").Value = 1054.6667
$ws.Range("N107").Value = -4865.2
$ws.Range("H113").Value = 3144.9412
$ws.Range("I113").Value = 2137.7
$ws.Range("K113").Value = 2137.7
$ws.Range("M113").Value = 32.30000000000018

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H123").Value = 66200
$ws.Range("J123").Value = 66200
$ws.Range("L123").Value = 66200
$ws.Range("N123").Value = -76000
$ws.Range("H132").Value = 25797.314
$ws.Range("I132").Value = 43395.383
$ws.Range("J132").Value = 9456.25
$ws.Range("K132").Value = 130186.149
$ws.Range("L132").Value = 28368.75
$ws.Range("M132").Value = -127656.149
$ws.Range("N132").Value = -33428.75
$ws.Range("H136").Value = 3535550.2
$ws.Range("I136").Value = 4620834.5
$ws.Range("J136").Value = 8377.5
$ws.Range("K136").Value = 13862503.5
$ws.Range("L136").Value = 25132.5
$ws.Range("M136").Value = -13859953.5
$ws.Range("N136").Value = -30232.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 179965
$ws.Range("J16").Value = 179965
$ws.Range("L16").Value = 179965
$ws.Range("N16").Value = -180549
$ws.Range("H27").Value = 169999.67
$ws.Range("J27").Value = 169997
$ws.Range("L27").Value = 169997
$ws.Range("N27").Value = -170135
$ws.Range("H62").Value = 9662.299999999999
$ws.Range("I62").Value = 9341.947
$ws.Range("K62").Value = 9341.947
$ws.Range("M62").Value = -8717.947
$ws.Range("H65").Value = 9662.299999999999
$ws.Range("I65").Value = 9341.947
$ws.Range("K65").Value = 46709.735
$ws.Range("M65").Value = -43589.735
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()
$ws.Range("H126").Value = 4013.5715
$ws.Range("I126").Value = 4819
$ws.Range("K126").Value = 14457
$ws.Range("M126").Value = -11987
$ws.Range("H132").Value = 3910.5557
$ws.Range("I132").Value = 2522.8474
$ws.Range("K132").Value = 7568.5422
$ws.Range("M132").Value = -5038.5422
$ws.Range("H137").Value = 147798.12
$ws.Range("J137").Value = 147798.12
$ws.Range("L137").Value = 147798.12
$ws.Range("N137").Value = -157998.12
